# cea/databases/CH/archetypes/construction_properties.xlsx
# "simplifying databases and fixing input locator"
#
# On the SUPPLY sheet (column G = type_dhw, column H = type_el):
#   - every row whose type_dhw was "T3" is changed to either "T19" or "T20"
#     (HOTEL and RETAIL building-use blocks get "T20", every other block
#     that had "T3" gets "T19")
#   - every data row's type_el ("T1") is changed to "T24"
# The SUPPLY sheet also becomes the active / selected sheet (with cell K224
# selected), while INDOOR_COMFORT stops being the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SUPPLY")

# --- column G (type_dhw): T3 -> T19 for most building-use blocks ---
$t19Range = $ws.Range("G24:G25,G39:G49,G63:G73,G75:G97,G99:G109,G111:G121,G123:G133,G146:G157,G170:G181,G195:G205,G207:G217,G219:G229")
foreach ($area in $t19Range.Areas) {
    $area.Value = "T19"
}

# --- column G (type_dhw): T3 -> T20 for the HOTEL and RETAIL blocks ---
$t20Range = $ws.Range("G27:G37,G51:G61")
foreach ($area in $t20Range.Areas) {
    $area.Value = "T20"
}

# --- column H (type_el): T1 -> T24 for every data row ---
$ws.Range("H2:H229").Value = "T24"

# Recalculate so the shared-formula cells in G183:G193 (which mirror
# G99:G109 for the LAB block) pick up the new cached values.
$excel.Calculate()

# Make SUPPLY the active sheet/tab and select cell K224, matching the
# workbook view recorded after the edit.
$ws.Activate() | Out-Null
$ws.Range("K224").Select() | Out-Null
